$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

# Row 2 (Q1)
$ws.Range("B2").Value = -0.2731398110005052
$ws.Range("C2").Value = 1.954904495510687
$ws.Range("D2").Value = 16.7369754246168
$ws.Range("E2").Value = 4.091084871353416
$ws.Range("F2").Value = 4.178015828355637
$ws.Range("G2").Value = 22

# Row 3 (Q2)
$ws.Range("B3").Value = -0.01488157163777362
$ws.Range("C3").Value = 1.752155569855608
$ws.Range("D3").Value = 10.88565705993317
$ws.Range("E3").Value = 3.29934191315983
$ws.Range("F3").Value = 3.380785024028695
$ws.Range("G3").Value = 21

# Row 4 (Q3)
$ws.Range("B4").Value = -0.6570990114281711
$ws.Range("C4").Value = 1.129967425434843
$ws.Range("D4").Value = 5.20882987611967
$ws.Range("E4").Value = 2.282286107419416
$ws.Range("F4").Value = 2.242426027718474
$ws.Range("G4").Value = 20

# Row 5 (Q4)
$ws.Range("B5").Value = -0.01990966418030321
$ws.Range("C5").Value = 0.7630979043556427
$ws.Range("D5").Value = 2.160275983970295
$ws.Range("E5").Value = 1.469787734324346
$ws.Range("F5").Value = 1.509924799371146
$ws.Range("G5").Value = 19

# Row 6 (Q5)
$ws.Range("B6").Value = -0.02640479156802359
$ws.Range("C6").Value = 0.7054922130995549
$ws.Range("D6").Value = 1.040078856359899
$ws.Range("E6").Value = 1.01984256449704
$ws.Range("F6").Value = 1.049057548473549
$ws.Range("G6").Value = 18

# Row 7 (Q6)
$ws.Range("B7").Value = -0.004386375610123202
$ws.Range("C7").Value = 0.6094268367774349
$ws.Range("D7").Value = 0.9227697331241523
$ws.Range("E7").Value = 0.960609042807818
$ws.Range("F7").Value = 0.9901628142054374
$ws.Range("G7").Value = 17

# Row 8 (final row)
$ws.Range("B8").Value = 0.1276921712949286
$ws.Range("C8").Value = 0.531544694086506
$ws.Range("D8").Value = 0.6230921363534133
$ws.Range("E8").Value = 0.7893618538752765
$ws.Range("F8").Value = 0.8045118408863187
$ws.Range("G8").Value = 16

$wb.Save()
